# Update "想去人数" (F column) figures across the relevant worksheets to
# reflect the latest scrape output (gh-pages generated data refresh).

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 630
$wsExhibit.Range("F3").Value = 476
$wsExhibit.Range("F7").Value = 38
$wsExhibit.Range("F8").Value = 1299
$wsExhibit.Range("F9").Value = 3955

# Sheet: 演出 (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 55
$wsShow.Range("F3").Value = 5

# Sheet: 全部类型 (All types, combined view)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 630
$wsAll.Range("F3").Value = 476
$wsAll.Range("F7").Value = 38
$wsAll.Range("F8").Value = 1299
$wsAll.Range("F9").Value = 3955
$wsAll.Range("F11").Value = 55
$wsAll.Range("F12").Value = 5
